$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Propagate the "Data analysis" row's direct formatting (font used
#    for B5 in the old layout) onto the three new "Data analysis" rows
#    (B11:B13) before we overwrite the old cells with new data.
# ---------------------------------------------------------------------
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B11:B13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Topic column (B) - reuses the existing topic strings, order does
#    not affect the shared-string table.
# ---------------------------------------------------------------------
$ws.Range("B2").Value2 = "Programming languages"
$ws.Range("B3").Value2 = "Programming languages"
$ws.Range("B4").Value2 = "Programming languages"
$ws.Range("B5").Value2 = "Machine learning"
$ws.Range("B6").Value2 = "Machine learning"
$ws.Range("B7").Value2 = "Machine learning"
$ws.Range("B8").Value2 = "Software engineering"
$ws.Range("B9").Value2 = "Software engineering"
$ws.Range("B10").Value2 = "Software engineering"
$ws.Range("B11").Value2 = "Data analysis"
$ws.Range("B12").Value2 = "Data analysis"
$ws.Range("B13").Value2 = "Data analysis"

# ---------------------------------------------------------------------
# 3. id column (A) and level column (D) - plain numbers, never touch
#    the shared-string table.
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = 1
$ws.Range("A3").Value2 = 2
$ws.Range("A4").Value2 = 3
$ws.Range("A5").Value2 = 4
$ws.Range("A6").Value2 = 5
$ws.Range("A7").Value2 = 6
$ws.Range("A8").Value2 = 7
$ws.Range("A9").Value2 = 8
$ws.Range("A10").Value2 = 9
$ws.Range("A11").Value2 = 10
$ws.Range("A12").Value2 = 11
$ws.Range("A13").Value2 = 12

$ws.Range("D2").Value2 = 4
$ws.Range("D3").Value2 = 3
$ws.Range("D4").Value2 = 2
$ws.Range("D5").Value2 = 4
$ws.Range("D6").Value2 = 3
$ws.Range("D7").Value2 = 3
$ws.Range("D8").Value2 = 2
$ws.Range("D9").Value2 = 3
$ws.Range("D10").Value2 = 3
$ws.Range("D11").Value2 = 4
$ws.Range("D12").Value2 = 4
$ws.Range("D13").Value2 = 3

# ---------------------------------------------------------------------
# 4. skills (C) and tooltip (E) columns - the order below reproduces
#    the original authoring order (skills for the first topic group,
#    then its tooltips, then the remaining skills which all reuse the
#    already-interned tooltip strings).
# ---------------------------------------------------------------------
$ws.Range("C2").Value2 = "python"
$ws.Range("C3").Value2 = "Java"
$ws.Range("C4").Value2 = "C"

$ws.Range("E2").Value2 = "helper.get_title_content"
$ws.Range("E3").Value2 = "1-2 years experience"
$ws.Range("E4").Value2 = "<1 year experience"

$ws.Range("C5").Value2 = "Supervised ML"
$ws.Range("C6").Value2 = "Unsupervised ML"
$ws.Range("C7").Value2 = "Deep Learning"
$ws.Range("C8").Value2 = "Front-end"
$ws.Range("C9").Value2 = "Databases"
$ws.Range("C10").Value2 = "Deployment"
$ws.Range("C11").Value2 = "Data Wrangling"
$ws.Range("C12").Value2 = "Data Visuzlization"
$ws.Range("C13").Value2 = "Statistics"

$ws.Range("E5").Value2 = "helper.get_title_content"
$ws.Range("E6").Value2 = "1-2 years experience"
$ws.Range("E7").Value2 = "<1 year experience"
$ws.Range("E8").Value2 = "helper.get_title_content"
$ws.Range("E9").Value2 = "1-2 years experience"
$ws.Range("E10").Value2 = "<1 year experience"
$ws.Range("E11").Value2 = "helper.get_title_content"
$ws.Range("E12").Value2 = "1-2 years experience"
$ws.Range("E13").Value2 = "<1 year experience"

# ---------------------------------------------------------------------
# 5. B5 no longer holds the "Data analysis" row, so drop the direct
#    formatting that used to live there (match the plain header style).
# ---------------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 6. Column C now needs to fit its new (longer) contents (mirrors the
#    author's "best fit" autosize after the new skill labels were added).
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15

# ---------------------------------------------------------------------
# 7. Move the active selection, matching where the author ended up.
# ---------------------------------------------------------------------
$ws.Range("I14").Select() | Out-Null

Write-Output "done"
